# Update Mongolia MSME summary sheet with refined (2-decimal) figures.
# The target cells hold their numbers as text (not numeric) values, so we
# force the Text number format on them first -- otherwise Excel would
# auto-convert a value like "22.23" into a genuine number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B12:D13").NumberFormat = "@"
$ws.Range("B15:D15").NumberFormat = "@"

$ws.Range("B12").Value = "22.23"
$ws.Range("C12").Value = "4.48"
$ws.Range("D12").Value = "26.72"

# Employment (% of total): Micro / SMEs / MSMEs
$ws.Range("B13").Value = "18.26"
$ws.Range("C13").Value = "33.82"
$ws.Range("D13").Value = "52.07"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B15").Value = "81.71"
$ws.Range("C15").Value = "16.47"
$ws.Range("D15").Value = "98.17"
